$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Good Morning" greeting text to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select cell E8, matching the selection recorded in the saved sheet view
$ws.Range("E8").Select()
